# Adapt column header formatting to respective input file names (#7)
#   Segmentname_old / ... -> Segmentname_FV2404 / ...
#   Segmentname_new / ... -> Segmentname_FV2410 / ...
# and turn the data range into an Excel Table ("Table1") with an
# autofilter + frozen header row, matching the regenerated merged AHB
# export.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- rename header row -------------------------------------------------
# columns A:J  -> "<name>_old" becomes "<name>_FV2404"
for ($c = 1; $c -le 10; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cell.Value = ($cell.Value2 -replace '_old$', '_FV2404')
}
# column K is "diff" - untouched
# columns L:U -> "<name>_new" becomes "<name>_FV2410"
for ($c = 12; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cell.Value = ($cell.Value2 -replace '_new$', '_FV2410')
}

# --- turn the used range into a proper Table1 (with autofilter) --------
$dataRange = $ws.Range("A1:U64")
$listObj = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$listObj.Name = "Table1"

# --- freeze the header row ---------------------------------------------
$ws.Activate()
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

# restore the cursor to A1 (matches the default sheet view selection)
$ws.Range("A1").Select() | Out-Null
